$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 1 ("100") expands into 13 single-value rows ---------------------
# Target values, in order: 0M, 0M, 0M, 20, 0.00002, 0.00005, 0.00003,
# 0.00001, 0.00003, 0.00003, 0.00004, 0.00068, 100.0
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"

$row1Values = @("0M","0M","20","0.00002","0.00005","0.00003","0.00001","0.00003","0.00003","0.00004","0.00068","100.0")
$anchorAfterRow1 = $t.Rows.Item(2)   # the unchanged "0" row stays a stable anchor
for ($i = $row1Values.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($anchorAfterRow1)
    $newRow.Cells.Item(1).Range.Text = $row1Values[$i]
}

# --- Row "70" (now at Item(15)) expands into 10 single-value rows --------
# Target values, in order: 0.00000 x9, 0.0
$row70Index = 15
$t.Rows.Item($row70Index).Cells.Item(1).Range.Text = "0.00000"

$row70Values = @("0.00000","0.00000","0.00000","0.00000","0.00000","0.00000","0.00000","0.00000","0.0")
$anchorAfterRow70 = $t.Rows.Item($row70Index + 1)   # the unchanged "0" row stays a stable anchor
for ($i = $row70Values.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($anchorAfterRow70)
    $newRow.Cells.Item(1).Range.Text = $row70Values[$i]
}

# --- Collapse the tab-separated "20 ... 100.0" row down to just "100" ----
$t.Rows.Item(45).Cells.Item(1).Range.Text = "100"

# --- Delete the now-orphaned empty row that followed it -------------------
$t.Rows.Item(46).Delete()

# --- Collapse the tab-separated "0 0 0 0 0 0 0 0 0 0" row down to "0" ----
$t.Rows.Item(46).Cells.Item(1).Range.Text = "0"

# --- Append the new trailing row containing "70" ---------------------------
$newLastRow = $t.Rows.Add()
$newLastRow.Cells.Item(1).Range.Text = "70"
